# Sprint Backlog update
# - Fills in the "Amount Remaining After..." columns (E, F, G) for the
#   sprint days that previously had no recorded value, and updates a
#   couple of existing estimates (D5, D15) to reflect real progress.
# - The weekly totals row (row 27, SUM formulas) and the burndown chart
#   series recalculate automatically from these cell values.
# - Updates the sheet's scroll position / selection to match where the
#   author left off editing (around row 13, cell F15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: Verify User Filtering and Implement Date Range Filtering ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0

# --- Row 4: Implement Cook user role ---
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# --- Row 5: Implement Stock item categories ---
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

# --- Row 6: Implement within inventory view remove food item ---
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

# --- Row 7: Implement stock changes log for removing food ---
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0

# --- Row 8: Implement Changes in Add Stock for picking date. ---
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# --- Row 9: Update Add Stock Page to Include Selection of Stock Category ---
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0

# --- Row 10: Implement Officer user role and add munitions item category ---
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

# --- Row 11: Implement Inventory View filtering for with officer role ---
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# --- Row 12: Implement View Inventory and View Model ---
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

# --- Row 13: Implement filtering view inventory with cook role ---
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# --- Row 14: Implement within inventory view remove munition item ---
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

# --- Row 15: Implement stock changes log for removing munition ---
$ws.Range("D15").Value = 0.5
$ws.Range("E15").Value = 0.5
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# --- Row 16: Implement food expiration report within view inventory for cook ---
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 0

# --- Row 17: Inform user about capacity remaining if compartment is overfilled. ---
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0

# --- Recalculate so the SUM(row 27) totals and the chart's cached values
#     pick up the new numbers ---
$excel.CalculateFullRebuild()

# --- Restore the view: scrolled to around row 13, with F15 selected ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F15").Select()
